$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "312.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.53%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.79%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.145"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.54%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07908"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.97%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.417"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.34%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.915"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.36%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.272"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.98%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.97%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9289"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.75%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1114"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-9.58%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1898"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.94%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09122"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.38%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03326"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-3.31%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09617"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.92%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001385"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.12%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005706"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-5.69%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.579"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.50%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3407"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.07%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.930"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "17.90%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.64%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2592"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.06%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04373"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.45%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004639"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "9.06%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.67%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003994"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02249"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.43%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05089"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.36%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007455"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-4.95%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009028"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-9.79%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1356"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.98%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.90%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008644"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.92%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006697"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.65%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003283"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2.22%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-40.72%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
